$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icosl"
$ws.Range("C2").Value = "Cd28"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.14972333333333
$ws.Range("H2").Value = 30.44917
$ws.Range("I2").Value = 0.4728485742924058
$ws.Range("J2").Value = 0.4728485742924059
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1791996666666667
$ws.Range("N2").Value = 0.537599
$ws.Range("O2").Value = 0.04251079199666429
$ws.Range("P2").Value = 0.04251079199666429
$ws.Range("Q2").Value = 1.818827038092222
$ws.Range("R2").Value = 16.36944334283
$ws.Range("S2").Value = 0.02010116738766372
$ws.Range("T2").Value = 0.02010116738766373

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Icosl"
$ws.Range("C3").Value = "Cd28"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.14972333333333
$ws.Range("H3").Value = 30.44917
$ws.Range("I3").Value = 0.4728485742924058
$ws.Range("J3").Value = 0.4728485742924059
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.036192666666667
$ws.Range("N3").Value = 12.108578
$ws.Range("O3").Value = 0.9574892080033358
$ws.Range("P3").Value = 0.9574892080033357
$ws.Range("Q3").Value = 40.96623888669556
$ws.Range("R3").Value = 368.69614998026
$ws.Range("S3").Value = 0.4527474069047421
$ws.Range("T3").Value = 0.4527474069047422

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Icosl"
$ws.Range("C4").Value = "Cd28"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.876575666666667
$ws.Range("H4").Value = 5.629727
$ws.Range("I4").Value = 0.08742466167732857
$ws.Range("J4").Value = 0.08742466167732858
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1791996666666667
$ws.Range("N4").Value = 0.537599
$ws.Range("O4").Value = 0.04251079199666429
$ws.Range("P4").Value = 0.04251079199666429
$ws.Range("Q4").Value = 0.3362817339414444
$ws.Range("R4").Value = 3.026535605473
$ws.Range("S4").Value = 0.003716491607943663
$ws.Range("T4").Value = 0.003716491607943663

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Icosl"
$ws.Range("C5").Value = "Cd28"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.876575666666667
$ws.Range("H5").Value = 5.629727
$ws.Range("I5").Value = 0.08742466167732857
$ws.Range("J5").Value = 0.08742466167732858
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.036192666666667
$ws.Range("N5").Value = 12.108578
$ws.Range("O5").Value = 0.9574892080033358
$ws.Range("P5").Value = 0.9574892080033357
$ws.Range("Q5").Value = 7.574220944245111
$ws.Range("R5").Value = 68.167988498206
$ws.Range("S5").Value = 0.08370817006938491
$ws.Range("T5").Value = 0.08370817006938491

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Icosl"
$ws.Range("C6").Value = "Cd28"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.840163333333335
$ws.Range("H6").Value = 26.52049
$ws.Range("I6").Value = 0.4118396621660297
$ws.Range("J6").Value = 0.4118396621660297
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1791996666666667
$ws.Range("N6").Value = 0.537599
$ws.Range("O6").Value = 0.04251079199666429
$ws.Range("P6").Value = 0.04251079199666429
$ws.Range("Q6").Value = 1.584154322612223
$ws.Range("R6").Value = 14.25738890351
$ws.Range("S6").Value = 0.01750763021431658
$ws.Range("T6").Value = 0.01750763021431658

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Icosl"
$ws.Range("C7").Value = "Cd28"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.840163333333335
$ws.Range("H7").Value = 26.52049
$ws.Range("I7").Value = 0.4118396621660297
$ws.Range("J7").Value = 0.4118396621660297
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.036192666666667
$ws.Range("N7").Value = 12.108578
$ws.Range("O7").Value = 0.9574892080033358
$ws.Range("P7").Value = 0.9574892080033357
$ws.Range("Q7").Value = 35.68060241813556
$ws.Range("R7").Value = 321.12542176322
$ws.Range("S7").Value = 0.3943320319517131
$ws.Range("T7").Value = 0.3943320319517132

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Icosl"
$ws.Range("C8").Value = "Cd28"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5985983333333333
$ws.Range("H8").Value = 1.795795
$ws.Range("I8").Value = 0.02788710186423574
$ws.Range("J8").Value = 0.02788710186423574
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1791996666666667
$ws.Range("N8").Value = 0.537599
$ws.Range("O8").Value = 0.04251079199666429
$ws.Range("P8").Value = 0.04251079199666429
$ws.Range("Q8").Value = 0.1072686218005556
$ws.Range("R8").Value = 0.9654175962050001
$ws.Range("S8").Value = 0.001185502786740314
$ws.Range("T8").Value = 0.001185502786740315

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Icosl"
$ws.Range("C9").Value = "Cd28"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5985983333333333
$ws.Range("H9").Value = 1.795795
$ws.Range("I9").Value = 0.02788710186423574
$ws.Range("J9").Value = 0.02788710186423574
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.036192666666667
$ws.Range("N9").Value = 12.108578
$ws.Range("O9").Value = 0.9574892080033358
$ws.Range("P9").Value = 0.9574892080033357
$ws.Range("Q9").Value = 2.416058203278889
$ws.Range("R9").Value = 21.74452382951
$ws.Range("S9").Value = 0.02670159907749543
$ws.Range("T9").Value = 0.02670159907749543
